$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Clear-ToEmptyText($range) {
    # Leaves the cell as an explicit empty text value (matching how this
    # sheet represents "blank" everywhere else) instead of deleting the
    # cell outright, and strips the quote-prefix style that Value="'"
    # would otherwise leave behind.
    $range.Value = "'"
    $range.Style = "Normal"
}

# --- Row 2: clear the duplicate registration (same competitor as row 5) ---
# enforce 1:1 competition <-> competitor relationship by wiping the
# duplicate row's data instead of leaving a second copy around.
Clear-ToEmptyText $ws.Range("A2")
Clear-ToEmptyText $ws.Range("B2")
Clear-ToEmptyText $ws.Range("C2")
Clear-ToEmptyText $ws.Range("D2")
Clear-ToEmptyText $ws.Range("H2")

# --- Row 6: Versenyengedelyszam was stored as text, fix it to a real number ---
$ws.Range("A6").Value = 6858

# --- Row 7: new competitor record; identifiers must stay text, not numbers ---
$a7 = $ws.Range("A7")
$a7.NumberFormat = "@"
$a7.Value = "6865"
$a7.Style = "Normal"

$ws.Range("B7").Value = "Gál László"
$ws.Range("C7").Value = "B.T.K. Szituációs Lövész és Szabadidős Sportegyesület"

# 1:1 verseny-versenyző kapcsolat kényszerítése: this competitor row no
# longer carries a direct Verseny_ID link.
Clear-ToEmptyText $ws.Range("V7")
